{"js": "// Update the \u98ce\u63a7\u603b\u76d1\u7b80\u4ecb (risk-control director bio) document:\n//   1. The leading \"\u57fa\u91d1\u98ce\u63a7\u603b\u76d1\uff1a\u66f9\u53d1\" label paragraph becomes empty.\n//   2. The biography paragraph is rewritten with updated career details\n//      (new job title, education, and a detailed work-history timeline).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length < 2) {\n  throw new Error(\"Expected at least 2 paragraphs, found \" + paragraphs.items.length);\n}\n\nconst labelPara = paragraphs.items[0];\nconst bioPara = paragraphs.items[1];\nlabelPara.load(\"text\");\nbioPara.load(\"text\");\nawait context.sync();\n\n// 1) Clear the label paragraph's text entirely (paragraph itself is kept).\nif (labelPara.text.length > 0) {\n  const labelMatches = labelPara.search(labelPara.text, { matchCase: true });\n  labelMatches.load(\"items\");\n  await context.sync();\n  labelMatches.items[0].delete();\n  await context.sync();\n}\n\n// 2) Replace the whole biography paragraph text with the revised biography.\nconst bioMatches = bioPara.search(bioPara.text, { matchCase: true });\nbioMatches.load(\"items\");\nawait context.sync();\n\nconst newBioText = \"\u66f9\u53d1\uff0c\u7537\uff0c1983\u5e742\u6708\u51fa\u751f\uff0c\u8eab\u4efd\u8bc1\u53f7370282198302205112\uff0c\u7533\u8bf7\u673a\u6784\u5408\u89c4\u98ce\u63a7\u603b\u76d1\u3002\u66f9\u53d1\u5148\u751f2008\u5e74\u6bd5\u4e1a\u4e8e\u534e\u4e2d\u79d1\u6280\u5927\u5b66\u5316\u5b66\u9662\uff0c\u5de5\u5b66\u7855\u58eb\uff0c\u540e\u5c31\u804c\u4e8e\u6b66\u6c49\u950b\u8fc5\u6295\u8d44\u7ba1\u7406\u6709\u9650\u516c\u53f8\uff0c\u4e8e2010\u5e74\uff0c\u901a\u8fc7\u300a\u8bc1\u5238\u57fa\u7840\u77e5\u8bc6\u300b\u548c\u300a\u6295\u8d44\u5206\u6790\u300b\u4e24\u95e8\u8d44\u683c\u8003\u8bd5\uff0c\u5e76\u83b7\u5f97\u8bc1\u5238\u4ece\u4e1a\u8d44\u683c\uff1b\u4e8e2017\u5e743\u6708\uff0c\u901a\u8fc7\u300a\u57fa\u91d1\u6cd5\u5f8b\u6cd5\u89c4\u3001\u804c\u4e1a\u9053\u5fb7\u4e0e\u4e1a\u52a1\u89c4\u8303\u300b\u548c\u300a\u79c1\u52df\u80a1\u6743\u6295\u8d44\u57fa\u91d1\u57fa\u7840\u77e5\u8bc6\u300b\uff0c\u5e76\u83b7\u5f97\u57fa\u91d1\u4ece\u4e1a\u8d44\u683c\u8bc1\u3002\u66f9\u53d1\u5148\u751f2008-2013\u5e74\u5c31\u804c\u6b66\u6c49\u950b\u8fc5\u6295\u8d44\u7ba1\u7406\u6709\u9650\u516c\u53f8\u671f\u95f4\uff0c\u5148\u540e\u62c5\u4efb\u4ea4\u6613\u5458\u3001\u98ce\u63a7\u4e13\u5458\u53ca\u7b56\u7565\u7814\u7a76\u5458\uff0c\u4ece\u4e8b\u6295\u8d44\u4ea4\u6613\u5e02\u573a\u5206\u6790\u7814\u7a76\u53ca\u98ce\u9669\u63a7\u5236\uff0c\u4e3b\u8981\u4f7f\u7528c++\u7b49\u7f16\u7a0b\u8bed\u8a00\u7f16\u5199CTA\u7b56\u7565\u53ca\u5957\u5229\u7b56\u7565\uff1b2013\u5e74\u81f32015\u5e74\u5c31\u804c\u4e8e\u667a\u6c47\u91cf\u5316\u6295\u8d44\u6709\u9650\u516c\u53f8\uff0c\u5148\u540e\u62c5\u4efb\u7b56\u7565\u7814\u7a76\u5458\u53ca\u98ce\u63a7\u603b\u76d1\u52a9\u7406\uff0c\u4e3b\u8981\u8d1f\u8d23\u7f16\u5199\u91cf\u5316\u6295\u8d44\u7b56\u7565\u3001\u5e02\u573a\u6570\u636e\u5206\u6790\u4ee5\u53ca\u534f\u52a9\u98ce\u63a7\u603b\u76d1\u5b8c\u6210\u5e02\u573a\u4ea4\u6613\u98ce\u63a7\u5efa\u6a21\u53ca\u5b9e\u65bd\uff1b\u671f\u95f4\u8d1f\u8d23\u4e24\u4e2a\u5957\u5229\u7b56\u7565\u7684\u8fd0\u7ef4\u53ca\u98ce\u9669\u7ba1\u7406\u63a7\u5236\uff0c\u7528python\u642d\u5efa\u4e86\u8986\u76d6\u6574\u4e2a\u516c\u53f8\u4ea4\u6613\u5e73\u53f0\u7684\u98ce\u63a7\u5e73\u53f0\uff0c\u5b9e\u65f6\u76d1\u6d4b\u8d44\u91d1\u98ce\u9669\u66b4\u9732\uff1b\u540c\u65f6\u548c\u56fd\u5185\u5f88\u591a\u4f18\u79c0\u7684\u79c1\u52df\u57fa\u91d1\u7ba1\u7406\u516c\u53f8\u5efa\u7acb\u4e86\u7d27\u5bc6\u7684\u5408\u4f5c\u5173\u7cfb\uff0c\u5e76\u5728\u5408\u4f5c\u65b9\u8f85\u52a9\u4e0b\u5f00\u5c55\u79c1\u52df\u6258\u7ba1\u4e1a\u52a1\u8425\u9500\u3001\u5ba2\u6237\u9002\u5f53\u6027\u7ba1\u7406\u3001\u5c3d\u804c\u8c03\u67e5\u7b49\u7814\u7a76\u5de5\u4f5c\uff1b2015\u5e74\u81f32016\u5e74\uff0c\u5c31\u804c\u4e8e\u9752\u5c9b\u7f8e\u6cf0\u79d1\u6280\u6709\u9650\u516c\u53f8\uff0c\u62c5\u4efb\u9ad8\u7ea7\u7814\u7a76\u5458\uff0c\u8d1f\u8d23\u7814\u53d1\u53ca\u751f\u4ea7\u6570\u636e\u5206\u6790\uff0c\u95ee\u9898\u5b9a\u4f4d\uff0c\u89e3\u51b3\u65b9\u6848\u53ca\u98ce\u9669\u6e90\u63a7\u5236\uff1b 2017\u5e743\u6708\u52a0\u5165\u7533\u8bf7\u673a\u6784\uff0c\u62c5\u4efb\u5408\u89c4\u98ce\u63a7\u603b\u76d1\u3002\";\nbioMatches.items[0].insertText(newBioText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Update the \u98ce\u63a7\u603b\u76d1\u7b80\u4ecb (risk-control director bio) document:\n#   1. The leading \"\u57fa\u91d1\u98ce\u63a7\u603b\u76d1\uff1a\u66f9\u53d1\" label paragraph becomes empty.\n#   2. The biography paragraph is rewritten with updated career details.\n\n$d = $word.ActiveDocument\n\n# 1) Clear the first paragraph's text (the \"\u57fa\u91d1\u98ce\u63a7\u603b\u76d1\uff1a\u66f9\u53d1\" label) via Find/Replace,\n#    leaving the (now empty) paragraph in place.\n$labelRange = $d.Paragraphs(1).Range\n$find = $labelRange.Find\n$find.ClearFormatting()\n$find.Text = \"\u57fa\u91d1\u98ce\u63a7\u603b\u76d1\uff1a\u66f9\u53d1\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2) Replace the whole biography paragraph (2nd paragraph) with the revised text.\n#    Assigning .Text directly (rather than Delete + InsertAfter) keeps the run's\n#    existing character formatting (e.g. the east-Asian font hint).\n$bioRange = $d.Paragraphs(2).Range\n$bioRange.End = $bioRange.End - 1\n$bioRange.Text = \"\u66f9\u53d1\uff0c\u7537\uff0c1983\u5e742\u6708\u51fa\u751f\uff0c\u8eab\u4efd\u8bc1\u53f7370282198302205112\uff0c\u7533\u8bf7\u673a\u6784\u5408\u89c4\u98ce\u63a7\u603b\u76d1\u3002\u66f9\u53d1\u5148\u751f2008\u5e74\u6bd5\u4e1a\u4e8e\u534e\u4e2d\u79d1\u6280\u5927\u5b66\u5316\u5b66\u9662\uff0c\u5de5\u5b66\u7855\u58eb\uff0c\u540e\u5c31\u804c\u4e8e\u6b66\u6c49\u950b\u8fc5\u6295\u8d44\u7ba1\u7406\u6709\u9650\u516c\u53f8\uff0c\u4e8e2010\u5e74\uff0c\u901a\u8fc7\u300a\u8bc1\u5238\u57fa\u7840\u77e5\u8bc6\u300b\u548c\u300a\u6295\u8d44\u5206\u6790\u300b\u4e24\u95e8\u8d44\u683c\u8003\u8bd5\uff0c\u5e76\u83b7\u5f97\u8bc1\u5238\u4ece\u4e1a\u8d44\u683c\uff1b\u4e8e2017\u5e743\u6708\uff0c\u901a\u8fc7\u300a\u57fa\u91d1\u6cd5\u5f8b\u6cd5\u89c4\u3001\u804c\u4e1a\u9053\u5fb7\u4e0e\u4e1a\u52a1\u89c4\u8303\u300b\u548c\u300a\u79c1\u52df\u80a1\u6743\u6295\u8d44\u57fa\u91d1\u57fa\u7840\u77e5\u8bc6\u300b\uff0c\u5e76\u83b7\u5f97\u57fa\u91d1\u4ece\u4e1a\u8d44\u683c\u8bc1\u3002\u66f9\u53d1\u5148\u751f2008-2013\u5e74\u5c31\u804c\u6b66\u6c49\u950b\u8fc5\u6295\u8d44\u7ba1\u7406\u6709\u9650\u516c\u53f8\u671f\u95f4\uff0c\u5148\u540e\u62c5\u4efb\u4ea4\u6613\u5458\u3001\u98ce\u63a7\u4e13\u5458\u53ca\u7b56\u7565\u7814\u7a76\u5458\uff0c\u4ece\u4e8b\u6295\u8d44\u4ea4\u6613\u5e02\u573a\u5206\u6790\u7814\u7a76\u53ca\u98ce\u9669\u63a7\u5236\uff0c\u4e3b\u8981\u4f7f\u7528c++\u7b49\u7f16\u7a0b\u8bed\u8a00\u7f16\u5199CTA\u7b56\u7565\u53ca\u5957\u5229\u7b56\u7565\uff1b2013\u5e74\u81f32015\u5e74\u5c31\u804c\u4e8e\u667a\u6c47\u91cf\u5316\u6295\u8d44\u6709\u9650\u516c\u53f8\uff0c\u5148\u540e\u62c5\u4efb\u7b56\u7565\u7814\u7a76\u5458\u53ca\u98ce\u63a7\u603b\u76d1\u52a9\u7406\uff0c\u4e3b\u8981\u8d1f\u8d23\u7f16\u5199\u91cf\u5316\u6295\u8d44\u7b56\u7565\u3001\u5e02\u573a\u6570\u636e\u5206\u6790\u4ee5\u53ca\u534f\u52a9\u98ce\u63a7\u603b\u76d1\u5b8c\u6210\u5e02\u573a\u4ea4\u6613\u98ce\u63a7\u5efa\u6a21\u53ca\u5b9e\u65bd\uff1b\u671f\u95f4\u8d1f\u8d23\u4e24\u4e2a\u5957\u5229\u7b56\u7565\u7684\u8fd0\u7ef4\u53ca\u98ce\u9669\u7ba1\u7406\u63a7\u5236\uff0c\u7528python\u642d\u5efa\u4e86\u8986\u76d6\u6574\u4e2a\u516c\u53f8\u4ea4\u6613\u5e73\u53f0\u7684\u98ce\u63a7\u5e73\u53f0\uff0c\u5b9e\u65f6\u76d1\u6d4b\u8d44\u91d1\u98ce\u9669\u66b4\u9732\uff1b\u540c\u65f6\u548c\u56fd\u5185\u5f88\u591a\u4f18\u79c0\u7684\u79c1\u52df\u57fa\u91d1\u7ba1\u7406\u516c\u53f8\u5efa\u7acb\u4e86\u7d27\u5bc6\u7684\u5408\u4f5c\u5173\u7cfb\uff0c\u5e76\u5728\u5408\u4f5c\u65b9\u8f85\u52a9\u4e0b\u5f00\u5c55\u79c1\u52df\u6258\u7ba1\u4e1a\u52a1\u8425\u9500\u3001\u5ba2\u6237\u9002\u5f53\u6027\u7ba1\u7406\u3001\u5c3d\u804c\u8c03\u67e5\u7b49\u7814\u7a76\u5de5\u4f5c\uff1b2015\u5e74\u81f32016\u5e74\uff0c\u5c31\u804c\u4e8e\u9752\u5c9b\u7f8e\u6cf0\u79d1\u6280\u6709\u9650\u516c\u53f8\uff0c\u62c5\u4efb\u9ad8\u7ea7\u7814\u7a76\u5458\uff0c\u8d1f\u8d23\u7814\u53d1\u53ca\u751f\u4ea7\u6570\u636e\u5206\u6790\uff0c\u95ee\u9898\u5b9a\u4f4d\uff0c\u89e3\u51b3\u65b9\u6848\u53ca\u98ce\u9669\u6e90\u63a7\u5236\uff1b 2017\u5e743\u6708\u52a0\u5165\u7533\u8bf7\u673a\u6784\uff0c\u62c5\u4efb\u5408\u89c4\u98ce\u63a7\u603b\u76d1\u3002\"\n"}
